$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the extra descriptive text from the Height row, keeping only "<n>cm".
$ws.Range("B3").Value = "175cm"
$ws.Range("C3").Value = "60cm"
$ws.Range("D3").Value = "170cm"
$ws.Range("E3").Value = "165cm"
$ws.Range("F3").Value = "181cm"
$ws.Range("G3").Value = "164cm"
$ws.Range("H3").Value = "152cm"
$ws.Range("I3").Value = "158cm"
$ws.Range("J3").Value = "178cm"
$ws.Range("K3").Value = "162cm"

# Column B width shrinks (best-fit recalculation) now that the text is shorter.
# (The host quantizes ColumnWidth to whole display pixels, so 10.5 is the
# input that lands closest to the recorded 11.4140625 best-fit width.)
$ws.Columns.Item(2).ColumnWidth = 10.5
